# Update "想去人数" (want-to-go count) figures on both the "展览" sheet
# and the "全部类型" sheet to reflect freshly regenerated data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet): row 3 and row 4 in column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 263
$wsExhibit.Range("F4").Value = 906

# Sheet "全部类型" (4th sheet): row 4 and row 5 in column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 263
$wsAll.Range("F5").Value = 906
